# Update the Model column values for all data rows (row 2 through 67)
# from the long publisher path to the short model name.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A2:A67").Value = "llama-3.1-8b-instruct-maas"

# Autofit / set width for column A to match bestFit width observed in the diff
$ws.Columns.Item(1).ColumnWidth = 47.6640625

# Adjust the view: scroll position and selection
$ws.Application.ActiveWindow.ScrollRow = 30
$ws.Range("A2:I67").Select()
